$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "92.244.19"
$ws.Range("E2").Value = "  +1.87%  "
$ws.Range("D3").Value = "3.111.78"
$ws.Range("E3").Value = "  -3.31%  "
$ws.Range("E4").Value = "  +0.07%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "237.58"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.98%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "615.36"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.73%  "
$ws.Range("E7").Value = "  -1.42%  "
$ws.Range("E8").Value = "  +4.98%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").Value = "3.108.73"
$ws.Range("E10").Value = "  -3.36%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.739"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("E12").Value = "  -1.16%  "
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("D14").Value = "92.133.07"
$ws.Range("E14").Value = "  -7.91%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "34.31"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -3.15%  "
$ws.Range("E16").Value = "  -2.28%  "
$ws.Range("D17").Value = "3.699.35"
$ws.Range("E17").Value = "  -3.04%  "
$ws.Range("D18").Value = "3.118.10"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("E19").Value = "  +0.68%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "14.66"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -3.75%  "
$ws.Range("E21").Value = "  -3.96%  "
$ws.Range("E22").Value = "  +1.37%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "446.72"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.17%  "
$ws.Range("E24").Value = "  -3.62%  "
$ws.Range("E25").Value = "  -4.12%  "
$ws.Range("E26").Value = "  -2.59%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "11.83"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.67%  "
$ws.Range("D28").Value = "3.276.25"
$ws.Range("E29").Value = "  -0.12%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.134"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -6.09%  "
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("E32").Value = "  -0.92%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "9.13"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -2.89%  "
$ws.Range("E34").Value = "  -0.69%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "7.88"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +2.41%  "
$ws.Range("E36").Value = "  -7.02%  "
$ws.Range("E37").Value = "  -3.87%  "
$ws.Range("E38").Value = "  -4.08%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "3.86"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +1.35%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "483.86"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -5.35%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.29"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -4.44%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "23.87"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +8.08%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.433"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -4.89%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "3.29"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -4.21%  "
$ws.Range("E45").Value = "  +0.05%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "162.07"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +3.56%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.89"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -2.30%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.693"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -6.37%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.39"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("E50").Value = "  +3.73%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "4.44"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -2.15%  "
